$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - same style as the rest of row 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-15
$data = @(
    @(2, 1, 3),
    @(3, 1, 6),
    @(4, 1, 5),
    @(5, 1, 5),
    @(6, 1, 4),
    @(7, 1, 1),
    @(8, 1, 6),
    @(9, 1, 7),
    @(10, 1, 6),
    @(11, 1, 5),
    @(12, 1, 6),
    @(13, 1, 5),
    @(14, 1, 4),
    @(15, 8, 9)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
